# Generate Report for Handback
# Adds a second handed-back file (9f75aa24-c862-4956-be99-e0355a2c7a3b.md) as a
# new row to each of the three tables (Overview, zh-cn, de-de), and refreshes
# the "Latest HO Xliff Generate Date" / handoff-handback timestamps for the
# first file (6bc5fd6b-83dd-44ab-a817-8de84405dc24.md, previously named
# d3e98c78-52c1-48c9-870d-bc0ad3910f2d.md).

$wb = $excel.ActiveWorkbook

$oldBase = "d3e98c78-52c1-48c9-870d-bc0ad3910f2d"
$newBase1 = "6bc5fd6b-83dd-44ab-a817-8de84405dc24"
$newBase2 = "9f75aa24-c862-4956-be99-e0355a2c7a3b"

$hash1 = "21724adc9a6908321f321651c8d25ffb706d7aaa"
$hash2 = "5f293120a6903d78b5fc3b11971a564f7003d7d4"

$dateGen        = "2016-09-03 03:05:40"
$zhHandoff      = "2016-09-03 03:05:36"
$zhHandback     = "2016-09-03 03:05:54"
$deHandoff      = "2016-09-03 03:05:40"
$deHandback     = "2016-09-03 03:06:04"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

# refresh the existing row for file 1 (rename + new generate date)
$wsOv.Range("A2").Value = "$newBase1.md"
$wsOv.Hyperlinks.Add($wsOv.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2794065f8bb41119f16ba75fa9639e2b8f7e7347/e2e/$newBase1.md", "", "", "e2e\$newBase1.md") | Out-Null
$wsOv.Range("G2").Value = $dateGen

# add a new row to the Overview table for file 2
$loOv = $wsOv.ListObjects.Item(1)
$loOv.ListRows.Add() | Out-Null

$wsOv.Range("A3").Value = "$newBase2.md"
$wsOv.Range("C3").Value = ".md"
$wsOv.Range("E3").Value = "Handed back: in sync with en-US"
$wsOv.Range("F3").Value = "Handed back: in sync with en-US"
$wsOv.Range("G3").Value = $dateGen
$wsOv.Hyperlinks.Add($wsOv.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2794065f8bb41119f16ba75fa9639e2b8f7e7347/e2e/$newBase2.md", "", "", "e2e\$newBase2.md") | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# refresh the existing row for file 1
$wsZh.Range("A2").Value = "$newBase1.md"
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2794065f8bb41119f16ba75fa9639e2b8f7e7347/e2e/$newBase1.md", "", "", "$newBase1.md") | Out-Null
$wsZh.Range("G2").Value = "$newBase1.$hash1.zh-cn.xlf"
$wsZh.Range("H2").Value = $zhHandoff
$wsZh.Range("I2").Value = "$newBase1.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f46c79755048fb4253baab7369a1e9ab37d3470d/e2e/$newBase1.md", "", "", "$newBase1.md") | Out-Null
$wsZh.Range("J2").Value = "$newBase1.$hash1.zh-cn.xlf"
$wsZh.Range("K2").Value = $zhHandback

# add a new row to the zh-cn table for file 2
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = "$newBase2.$hash2.zh-cn.xlf"
$wsZh.Range("H3").Value = $zhHandoff
$wsZh.Range("J3").Value = "$newBase2.$hash2.zh-cn.xlf"
$wsZh.Range("K3").Value = $zhHandback
$wsZh.Range("M3").Value = "True"
$wsZh.Range("O3").Value = "False"

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2794065f8bb41119f16ba75fa9639e2b8f7e7347/e2e/$newBase2.md", "", "", "$newBase2.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f46c79755048fb4253baab7369a1e9ab37d3470d/e2e/$newBase2.md", "", "", "$newBase2.md") | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# refresh the existing row for file 1
$wsDe.Range("A2").Value = "$newBase1.md"
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2794065f8bb41119f16ba75fa9639e2b8f7e7347/e2e/$newBase1.md", "", "", "$newBase1.md") | Out-Null
$wsDe.Range("G2").Value = "$newBase1.$hash1.de-de.xlf"
$wsDe.Range("I2").Value = "$newBase1.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/565cf2ab56a943b0e2fa97f919f674667a072cbc/e2e/$newBase1.md", "", "", "$newBase1.md") | Out-Null
$wsDe.Range("J2").Value = "$newBase1.$hash1.de-de.xlf"
$wsDe.Range("K2").Value = $deHandback

# add a new row to the de-de table for file 2
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = "$newBase2.$hash2.de-de.xlf"
$wsDe.Range("H3").Value = $deHandoff
$wsDe.Range("J3").Value = "$newBase2.$hash2.de-de.xlf"
$wsDe.Range("K3").Value = $deHandback
$wsDe.Range("M3").Value = "True"
$wsDe.Range("O3").Value = "False"

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2794065f8bb41119f16ba75fa9639e2b8f7e7347/e2e/$newBase2.md", "", "", "$newBase2.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/565cf2ab56a943b0e2fa97f919f674667a072cbc/e2e/$newBase2.md", "", "", "$newBase2.md") | Out-Null

# ---------------------------------------------------------------------------
# Apply consistent number formatting to the datetime columns for the new rows
# ---------------------------------------------------------------------------
$wsOv.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

Write-Host "Handback report rows added."
